$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: staging / prod
$ws.Range("A22").Formula = '=IF(D22="prod",B22&".bcparks.ca",D22&"-"&B22&".bcparks.ca")'
$ws.Range("B22").Value = "staging"
$ws.Range("C22").Value = "Staging/Gatsby"
$ws.Range("D22").Value = "prod"

# Row 23: staging / dev
$ws.Range("A23").Formula = '=IF(D23="prod",B23&".bcparks.ca",D23&"-"&B23&".bcparks.ca")'
$ws.Range("B23").Value = "staging"
$ws.Range("C23").Value = "Staging/Gatsby"
$ws.Range("D23").Value = "dev"

# Row 24: staging / test
$ws.Range("A24").Formula = '=IF(D24="prod",B24&".bcparks.ca",D24&"-"&B24&".bcparks.ca")'
$ws.Range("B24").Value = "staging"
$ws.Range("C24").Value = "Staging/Gatsby"
$ws.Range("D24").Value = "test"

# Row 25: staging / tools
$ws.Range("A25").Formula = '=IF(D25="prod",B25&".bcparks.ca",D25&"-"&B25&".bcparks.ca")'
$ws.Range("B25").Value = "staging"
$ws.Range("C25").Value = "Staging/Gatsby"
$ws.Range("D25").Value = "tools"

$ws.Range("A22:A25").Select()
